$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the "v","l","s","q","c" columns (currently C:G) two columns to the
# left (onto A:E), overwriting the now-obsolete "nr"/"nm" columns. This
# engine's multi-column Range.Cut(Destination) has a data-loss bug, so the
# move is emulated with Copy + Clear of the vacated source cells.
$src = $ws.Range("C1:G7")
$src.Copy($ws.Range("A1"))
$ws.Range("F1:G7").Clear()

# The data-validation rule originally watched A1:E1. Two of the three
# surviving watched cells (C1->A1, D1->B1, E1->C1) and (A1, B1) now all
# collapse onto A1:C1, so rebuild the rule on that range with its settings
# unchanged.
$formula1 = 'AND(ISNUMBER(SUMPRODUCT(SEARCH(MID(A1,ROW(INDIRECT("1:"&LEN(A1))),1),"0123456789abcdefghijklmnopqrstuvwxyzABCDEFGHIJKLMNOPQRSTUVWXYZ_"))),ISNUMBER(SEARCH(LEFT(A1,1),"abcdefghijklmnopqrstuvwxyzABCDEFGHIJKLMNOPQRSTUVWXYZ_")),NOT(ISNUMBER(SEARCH("~*",A1))))'
$ws.Range("A1:E1").Validation.Delete()
$validation = $ws.Range("A1:C1").Validation
$validation.Add(7, 1, 1, $formula1)
$validation.ErrorTitle = "Invalid variable name"
$validation.ErrorMessage = "Variable names need to start with a letter or an underscore (_), followed by only letters, numbers and underscores."
$validation.IgnoreBlank = $true
$validation.ShowInput = $true
$validation.ShowError = $true

# Leave the cursor where the author's session ended up.
$ws.Range("G6").Select() | Out-Null
